$wb = $excel.ActiveWorkbook

# ---- Update Timestamp (shared across all rows in FBS!AK column) ----
$wsFbs = $wb.Worksheets.Item("FBS")
$wsFbs.Range("AK2:AK57").Value = "2024-10-16T10:01:35.882115"

# ---- FBS sheet data updates ----
$ws = $wsFbs
    $ws.Range("N5").Value = "W"
    $ws.Range("Q5").Value = "W"
    $ws.Range("AB6").Value = -2.5
    $ws.Range("AF6").Value = 0
    $ws.Range("M6").Value = "WSW"
    $ws.Range("N6").Value = "WSW"
    $ws.Range("O6").Value = 62.18000000000001
    $ws.Range("P6").Value = 6.8
    $ws.Range("Q6").Value = "WSW"
    $ws.Range("U6").Value = 1
    $ws.Range("AE7").Value = -0.07766990291262135
    $ws.Range("N7").Value = "W"
    $ws.Range("O7").Value = 69.2
    $ws.Range("P7").Value = 11.7
    $ws.Range("U7").Value = 3.7
    $ws.Range("Y7").Value = 47.5
    $ws.Range("Z7").Value = -110
    $ws.Range("O8").Value = 49.58000000000001
    $ws.Range("P8").Value = 4
    $ws.Range("U8").Value = 0.1
    $ws.Range("N9").Value = "SE"
    $ws.Range("O9").Value = 43.85
    $ws.Range("P9").Value = 5.1
    $ws.Range("Q9").Value = "SSE"
    $ws.Range("U9").Value = 1.1
    $ws.Range("M13").Value = "S"
    $ws.Range("M14").Value = "SE"
    $ws.Range("M15").Value = "SE"
    $ws.Range("N15").Value = "SE"
    $ws.Range("Q15").Value = "SE"
    $ws.Range("M16").Value = "W"
    $ws.Range("N16").Value = "WNW"
    $ws.Range("Q16").Value = "W"
    $ws.Range("M18").Value = "W"
    $ws.Range("N18").Value = "W"
    $ws.Range("Q18").Value = "W"
    $ws.Range("AE19").Value = -0.01904761904761905
    $ws.Range("Y19").Value = 51.5
    $ws.Range("N20").Value = "W"
    $ws.Range("Q20").Value = "W"
    $ws.Range("Z22").Value = -112
    $ws.Range("N24").Value = "S"
    $ws.Range("Q24").Value = "S"
    $ws.Range("M37").Value = "S"
    $ws.Range("Q37").Value = "SW"
    $ws.Range("M38").Value = "SE"
    $ws.Range("N38").Value = "SE"
    $ws.Range("Q38").Value = "SE"
    $ws.Range("M39").Value = "SE"
    $ws.Range("N39").Value = "SE"
    $ws.Range("Q39").Value = "SE"
    $ws.Range("M40").Value = "SSE"
    $ws.Range("N40").Value = "SE"
    $ws.Range("Q40").Value = "SSE"
    $ws.Range("M41").Value = "W"
    $ws.Range("Q45").Value = "NNE"
    $ws.Range("N46").Value = "WNW"
    $ws.Range("Q46").Value = "WNW"
    $ws.Range("AB48").Value = -4
    $ws.Range("AF48").Value = 0.5
    $ws.Range("M48").Value = "WNW"
    $ws.Range("M52").Value = "SSE"
    $ws.Range("M54").Value = "SE"
    $ws.Range("N54").Value = "S"
    $ws.Range("Q54").Value = "SSE"
    $ws.Range("M56").Value = "SE"
    $ws.Range("N56").Value = "SE"
    $ws.Range("Q56").Value = "SE"
    $ws.Range("N57").Value = "WNW"
    $ws.Range("Q57").Value = "W"

# ---- Other sheet data updates ----
$ws = $wb.Worksheets.Item("Other")
    $ws.Range("O3").Value = "WNW"
    $ws.Range("P3").Value = "WNW"
    $ws.Range("S3").Value = "WNW"
    $ws.Range("O4").Value = "WNW"
    $ws.Range("P4").Value = "WNW"
    $ws.Range("S4").Value = "WNW"
    $ws.Range("S5").Value = "SW"
    $ws.Range("P10").Value = "S"
    $ws.Range("S10").Value = "S"
    $ws.Range("O11").Value = "SE"
    $ws.Range("P11").Value = "SE"
    $ws.Range("S11").Value = "SE"
    $ws.Range("O15").Value = "SSE"
    $ws.Range("P15").Value = "SE"
    $ws.Range("S15").Value = "SSE"
    $ws.Range("O23").Value = "W"
    $ws.Range("P23").Value = "W"
    $ws.Range("S23").Value = "W"
    $ws.Range("P25").Value = "WNW"
    $ws.Range("S25").Value = "WNW"
    $ws.Range("O26").Value = "W"
    $ws.Range("P26").Value = "W"
    $ws.Range("S26").Value = "W"
    $ws.Range("O27").Value = "W"
    $ws.Range("O28").Value = "S"
    $ws.Range("S28").Value = "S"
    $ws.Range("O29").Value = "SSE"
    $ws.Range("P29").Value = "S"
    $ws.Range("S29").Value = "SSE"
    $ws.Range("O32").Value = "S"
    $ws.Range("P37").Value = "W"
    $ws.Range("S37").Value = "W"
    $ws.Range("O38").Value = "S"
    $ws.Range("S38").Value = "S"

